$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    [void]$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-23 Monday" "2024-09-24 Tuesday"

Replace-Text "70×67=4690" "58×42=2436"
Replace-Text "48×65=3120" "46×89=4094"
Replace-Text "35×90=3150" "59×99=5841"
Replace-Text "69×81=5589" "50×52=2600"
Replace-Text "48×50=2400" "79×69=5451"

Replace-Text "81×56=4536" "68×35=2380"
Replace-Text "33×22=726" "60×14=840"
Replace-Text "76×17=1292" "76×56=4256"
Replace-Text "18×99=1782" "35×45=1575"

Replace-Text "90×74=6660" "88×30=2640"
Replace-Text "39×69=2691" "21×57=1197"
Replace-Text "40×78=3120" "77×30=2310"
Replace-Text "56×23=1288" "99×49=4851"
Replace-Text "91×26=2366" "63×65=4095"

Replace-Text "53×87=4611" "29×13=377"
Replace-Text "43×14=602" "35×74=2590"
Replace-Text "47×35=1645" "43×23=989"
Replace-Text "28×25=700" "33×99=3267"

Replace-Text "22×77=1694" "53×64=3392"
Replace-Text "28×99=2772" "30×62=1860"
Replace-Text "40×16=640" "44×26=1144"
Replace-Text "53×69=3657" "82×80=6560"
Replace-Text "99×59=5841" "43×37=1591"

# Handle the two distinct occurrences of "37×70=2590" via direct table cell access
# Row 5 (1-based), column 2 -> 17×67=1139
# Row 15 (1-based), column 4 -> 50×35=1750
$table = $d.Tables.Item(1)
$table.Cell(5, 2).Range.Text = "17×67=1139"
$table.Cell(15, 4).Range.Text = "50×35=1750"

Write-Host "Done"
